$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete createPlayerCommandFactory with undo/redo:
# "1. Create Player" row now fully supports undo/redo -> mark as "Y (OK)"
# and fill in the "command pattern" / "command factory pattern" columns.
$ws.Range("C3").Value = "Y (OK)"
$ws.Range("E3").Value = $ws.Range("D3").Value()

# "9. show undo/ redo list" and "10. undo" rows also gain command-pattern support.
$ws.Range("D12").Value = $ws.Range("D3").Value()
$ws.Range("D13").Value = $ws.Range("D3").Value()

# The "Implement without command pattern" column is no longer needed for review, hide it.
$ws.Columns.Item(2).Hidden = $true

# Move the active selection to C3.
[void]$ws.Range("C3").Select()
